$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-5: date (D), volumen (M), precio minimo (N),
# precio maximo (O), precio promedio ponderado (P), precio $/Kg (S)
$updates = @(
    @{ Row = 2; D = 44971; M = 25; N = 28000; O = 28000; P = 28000; S = 3500 },
    @{ Row = 3; D = 44981; M = 30; N = 25000; O = 25000; P = 25000; S = 3125 },
    @{ Row = 4; D = 44973; M = 55; N = 28000; O = 28000; P = 28000; S = 3500 },
    @{ Row = 5; D = 44980; M = 50; N = 25000; O = 25000; P = 25000; S = 3125 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("M$r").Value = $u.M
    $ws.Range("N$r").Value = $u.N
    $ws.Range("O$r").Value = $u.O
    $ws.Range("P$r").Value = $u.P
    $ws.Range("S$r").Value = $u.S
}
